$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B4").Value = 7.966999999999999
$ws.Range("A10").Value = -21.776
$ws.Range("A12").Value = -21.53
$ws.Range("B12").Value = 5.397
$ws.Range("B17").Value = 5.179
$ws.Range("A18").Value = -21.879
$ws.Range("B26").Value = 5.720000000000001
$ws.Range("B27").Value = 6.222
$ws.Range("B28").Value = 5.739999999999999
$ws.Range("A37").Value = -21.032
$ws.Range("B37").Value = 7.816999999999998
$ws.Range("A55").Value = -22.184
$ws.Range("B65").Value = 5.912999999999999
$ws.Range("A68").Value = -21.567
$ws.Range("B73").Value = 6.528
$ws.Range("A77").Value = -20.93
$ws.Range("A78").Value = -20.582
$ws.Range("B84").Value = 5.863000000000001
$ws.Range("B85").Value = 5.964
$ws.Range("B93").Value = 5.587
$ws.Range("B95").Value = 5.898999999999999
$ws.Range("B98").Value = 6.423
$ws.Range("B99").Value = 5.414999999999999
$ws.Range("B101").Value = 5.372

$wb.Save()
